$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.161.33"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.516.40"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.36"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.19"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.550"
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.16"
$ws.Range("E10").Value = "  +4.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.28"
$ws.Range("E11").Value = "  +11.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0820"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.24"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.912.53"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.517.72"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.851"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.011.22"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.22"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0944"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.70"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.95"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "274.99"
$ws.Range("E24").Value = "  +12.23%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.96"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("E30").Value = "  +6.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.50"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.57"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.69"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0784"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.49"
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.031.44"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.14"
$ws.Range("E46").Value = "  +4.64%  "
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  +6.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.03"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.18"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.16"
$ws.Range("E51").Value = "  +4.36%  "
